$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9192037738362728
$ws.Range("C2").Value = 0.3034642044029852
$ws.Range("D2").Value = 0.0789843654385578
$ws.Range("E2").Value = 0.09074448073276642
$ws.Range("G2").Value = 0.8581293149973135
$ws.Range("H2").Value = 0.8286167870744805
$ws.Range("I2").Value = 0.6040612733691333
$ws.Range("M2").Value = 0.3909843463112708
$ws.Range("N2").Value = 1.115706564159545
$ws.Range("B3").Value = 0.8148995754344242
$ws.Range("C3").Value = 0.2649955500315855
$ws.Range("D3").Value = 0.07156238137076798
$ws.Range("E3").Value = 0.08585969797372783
$ws.Range("G3").Value = 0.8251535960428669
$ws.Range("H3").Value = 0.819283179797452
$ws.Range("I3").Value = 0.5984475280528017
$ws.Range("M3").Value = 0.351131231888715
$ws.Range("N3").Value = 1.132679338412453
$ws.Range("B4").Value = 0.7510506840767448
$ws.Range("C4").Value = 0.2413982762189448
$ws.Range("D4").Value = 0.06704538856395459
$ws.Range("E4").Value = 0.08293268238793416
$ws.Range("G4").Value = 0.8056340871627583
$ws.Range("H4").Value = 0.814119184211151
$ws.Range("I4").Value = 0.5954659778289368
$ws.Range("M4").Value = 0.3268221343125433
$ws.Range("N4").Value = 1.143622773648699
$ws.Range("B5").Value = 0.725080035947343
$ws.Range("C5").Value = 0.2317873506584931
$ws.Range("D5").Value = 0.06521463384248705
$ws.Range("E5").Value = 0.08175782585947644
$ws.Range("G5").Value = 0.7978610902309384
$ws.Range("H5").Value = 0.8121567391399083
$ws.Range("I5").Value = 0.5943673224641302
$ws.Range("M5").Value = 0.3169557691720541
$ws.Range("N5").Value = 1.148213243699061
$ws.Range("B6").Value = 0.7207705303457885
$ws.Range("C6").Value = 0.2301917596261092
$ws.Range("D6").Value = 0.06491123582472369
$ws.Range("E6").Value = 0.08156381741032348
$ws.Range("G6").Value = 0.7965812963822003
$ws.Range("H6").Value = 0.8118394305222836
$ws.Range("I6").Value = 0.5941918997644677
$ws.Range("M6").Value = 0.3153198507294235
$ws.Range("N6").Value = 1.148983385507606
$ws.Range("B7").Value = 0.7507002400291753
$ws.Range("C7").Value = 0.2412686397418895
$ws.Range("D7").Value = 0.06702065822840098
$ws.Range("E7").Value = 0.08291676560084227
$ws.Range("G7").Value = 0.8055285256353386
$ws.Range("H7").Value = 0.8140921441710702
$ws.Range("I7").Value = 0.5954506907525641
$ws.Range("M7").Value = 0.326688912823947
$ws.Range("N7").Value = 1.143684152671735
$ws.Range("B8").Value = 0.8831989790732564
$ws.Range("C8").Value = 0.2901951782383492
$ws.Range("D8").Value = 0.07641685083952154
$ws.Range("E8").Value = 0.08904506064516937
$ws.Range("G8").Value = 0.8466072547403201
$ws.Range("H8").Value = 0.8252805235282494
$ws.Range("I8").Value = 0.6020286790651781
$ws.Range("M8").Value = 0.3772091885654518
$ws.Range("N8").Value = 1.121450208322694
$ws.Range("B9").Value = 1.144610042185093
$ws.Range("C9").Value = 0.3863486459703722
$ws.Range("D9").Value = 0.09516833612130426
$ws.Range("E9").Value = 0.1016477395342221
$ws.Range("G9").Value = 0.9330154245565154
$ws.Range("H9").Value = 0.8517492849637733
$ws.Range("I9").Value = 0.6186513605328656
$ws.Range("M9").Value = 0.4775910859238763
$ws.Range("N9").Value = 1.082007065851856
$ws.Range("B10").Value = 1.337708169686266
$ws.Range("C10").Value = 0.4571670742202514
$ws.Range("D10").Value = 0.1091548333229753
$ws.Range("E10").Value = 0.1112809510768997
$ws.Range("G10").Value = 1.000184415159566
$ws.Range("H10").Value = 0.8740022600259181
$ws.Range("I10").Value = 0.6331794395959207
$ws.Range("M10").Value = 0.5521979964106265
$ws.Range("N10").Value = 1.055582499845588
$ws.Range("B11").Value = 1.425796951753512
$ws.Range("C11").Value = 0.4894328082815491
$ws.Range("D11").Value = 0.1155657445505653
$ws.Range("E11").Value = 0.1157483972780753
$ws.Range("G11").Value = 1.031567299499017
$ws.Range("H11").Value = 0.8847452899994437
$ws.Range("I11").Value = 0.6403014630517774
$ws.Range("M11").Value = 0.586337358429688
$ws.Range("N11").Value = 1.044120241954776
$ws.Range("B12").Value = 1.459190501041064
$ws.Range("C12").Value = 0.501658932514772
$ws.Range("D12").Value = 0.1180005245195304
$ws.Range("E12").Value = 0.1174526562463285
$ws.Range("G12").Value = 1.043572118689696
$ws.Range("H12").Value = 0.8889033318841371
$ws.Range("I12").Value = 0.6430729588099666
$ws.Range("M12").Value = 0.5992947394345691
$ws.Range("N12").Value = 1.039860485691165
$ws.Range("B13").Value = 1.451996980326044
$ws.Range("C13").Value = 0.4990254619008283
$ws.Range("D13").Value = 0.1174758326027785
$ws.Range("E13").Value = 0.117085051712543
$ws.Range("G13").Value = 1.040981264106648
$ws.Range("H13").Value = 0.8880038135993118
$ws.Range("I13").Value = 0.642472740350776
$ws.Range("M13").Value = 0.5965028096115077
$ws.Range("N13").Value = 1.040774301646863
$ws.Range("B14").Value = 1.428543529042088
$ws.Range("C14").Value = 0.4904384993213284
$ws.Range("D14").Value = 0.1157659120656547
$ws.Range("E14").Value = 0.1158883549000009
$ws.Range("G14").Value = 1.032552511621532
$ws.Range("H14").Value = 0.8850855688809816
$ws.Range("I14").Value = 0.6405279774434831
$ws.Range("M14").Value = 0.5874027735338672
$ws.Range("N14").Value = 1.043768168229292
$ws.Range("B15").Value = 1.414182352111425
$ws.Range("C15").Value = 0.4851797690985791
$ws.Range("D15").Value = 0.1147194663761155
$ws.Range("E15").Value = 0.1151569840754973
$ws.Range("G15").Value = 1.027405443846817
$ws.Range("H15").Value = 0.8833097877864589
$ws.Range("I15").Value = 0.6393464837948315
$ws.Range("M15").Value = 0.5818326086043584
$ws.Range("N15").Value = 1.045612527839097
$ws.Range("B16").Value = 1.331956395721988
$ws.Range("C16").Value = 0.4550594904988543
$ws.Range("D16").Value = 0.108736851539291
$ws.Range("E16").Value = 0.1109907318600776
$ws.Range("G16").Value = 0.9981502744867612
$ws.Range("H16").Value = 0.8733127166123893
$ws.Range("I16").Value = 0.6327243882448528
$ws.Range("M16").Value = 0.549971002609638
$ws.Range("N16").Value = 1.056342860571281
$ws.Range("B17").Value = 1.281577276372275
$ws.Range("C17").Value = 0.4365949176651611
$ws.Range("D17").Value = 0.1050792091240851
$ws.Range("E17").Value = 0.1084569004234339
$ws.Range("G17").Value = 0.9804163801074139
$ws.Range("H17").Value = 0.8673391598455282
$ws.Range("I17").Value = 0.628793904761082
$ws.Range("M17").Value = 0.5304767663053127
$ws.Range("N17").Value = 1.063068914464669
$ws.Range("B18").Value = 1.252623792523025
$ws.Range("C18").Value = 0.4259792612343176
$ws.Range("D18").Value = 0.1029799749822899
$ws.Range("E18").Value = 0.107007515077008
$ws.Range("G18").Value = 0.9702940791732715
$ws.Range("H18").Value = 0.8639616449789855
$ws.Range("I18").Value = 0.6265814471496043
$ws.Range("M18").Value = 0.5192830056827091
$ws.Range("N18").Value = 1.066990095373878
$ws.Range("B19").Value = 1.24282461461604
$ws.Range("C19").Value = 0.4223857679136245
$ws.Range("D19").Value = 0.1022699864594188
$ws.Range("E19").Value = 0.1065181450623172
$ws.Range("G19").Value = 0.9668801430184999
$ws.Range("H19").Value = 0.862828070420278
$ws.Range("I19").Value = 0.6258406116075506
$ws.Range("M19").Value = 0.5154961986340822
$ws.Range("N19").Value = 1.06832674867864
$ws.Range("B20").Value = 1.286937808353912
$ws.Range("C20").Value = 0.4385600146027286
$ws.Range("D20").Value = 0.1054680999276201
$ws.Range("E20").Value = 0.1087258003764404
$ws.Range("G20").Value = 0.9822961217801378
$ws.Range("H20").Value = 0.8679690143495407
$ws.Range("I20").Value = 0.6292073120701076
$ws.Range("M20").Value = 0.532550009082982
$ws.Range("N20").Value = 1.062347474388535
$ws.Range("B21").Value = 1.435431387366009
$ws.Range("C21").Value = 0.4929604821854809
$ws.Range("D21").Value = 0.1162679630931649
$ws.Range("E21").Value = 0.1162395113821972
$ws.Range("G21").Value = 1.035024947365883
$ws.Range("H21").Value = 0.8859402826628298
$ws.Range("I21").Value = 0.6410971724178225
$ws.Range("M21").Value = 0.5900748668885996
$ws.Range("N21").Value = 1.042886601802635
$ws.Range("B22").Value = 1.532692074913371
$ws.Range("C22").Value = 0.5285600253384928
$ws.Range("D22").Value = 0.1233677938461852
$ws.Range("E22").Value = 0.1212233439705344
$ws.Range("G22").Value = 1.070191269799636
$ws.Range("H22").Value = 0.8982098355002393
$ws.Range("I22").Value = 0.6493027146487904
$ws.Range("M22").Value = 0.627843253693058
$ws.Range("N22").Value = 1.030638697033595
$ws.Range("B23").Value = 1.480762637146142
$ws.Range("C23").Value = 0.5095555007423513
$ws.Range("D23").Value = 0.1195746338030403
$ws.Range("E23").Value = 0.1185565902196899
$ws.Range("G23").Value = 1.051357227201862
$ws.Range("H23").Value = 0.8916131284816231
$ws.Range("I23").Value = 0.6448832245247118
$ws.Range("M23").Value = 0.6076695167069488
$ws.Range("N23").Value = 1.037132395104564
$ws.Range("B24").Value = 1.284514280306439
$ws.Range("C24").Value = 0.4376715947265666
$ws.Range("D24").Value = 0.1052922712162285
$ws.Range("E24").Value = 0.1086042078419567
$ws.Range("G24").Value = 0.9814460628529957
$ws.Range("H24").Value = 0.86768408036869
$ws.Range("I24").Value = 0.6290202635642643
$ws.Range("M24").Value = 0.5316126533016501
$ws.Range("N24").Value = 1.06267346838758
$ws.Range("B25").Value = 1.073714411993649
$ws.Range("C25").Value = 0.3603093929914394
$ws.Range("D25").Value = 0.09005954510001857
$ws.Range("E25").Value = 0.09817396096432418
$ws.Range("G25").Value = 0.9090011150524901
$ws.Range("H25").Value = 0.8440992916344499
$ws.Range("I25").Value = 0.6137511959724549
$ws.Range("M25").Value = 0.450288873531079
$ws.Range("N25").Value = 1.092230371868256
